$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the bookmark that currently sits at the end of the
#    ":set tabstop=4" paragraph (paragraph 2).  It will be re-created later
#    at the new cursor position (the last empty list paragraph).
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 2. The first empty bulleted placeholder paragraph (originally paragraph 5)
#    is where the new "Beyond compare 4 license" block gets typed in.  We
#    insert 16 extra empty paragraphs in front of it (it + the following
#    empty bulleted paragraph give us the 18 slots needed for the license
#    block, the new blank/bookmark line, and the still-bulleted blank
#    line that used to be paragraph 6).
# ---------------------------------------------------------------------------
$placeholder = $d.Paragraphs.Item(5)
$anchor = $placeholder.Range.Duplicate
$anchor.Collapse(1)   # wdCollapseStart

for ($i = 0; $i -lt 16; $i++) {
    $anchor.InsertParagraphBefore()
}

# ---------------------------------------------------------------------------
# 3. Fill in the text of the 16 new paragraphs (indices 5-20) plus the
#    original placeholder (now at index 5 as well, since InsertParagraphBefore
#    pushed everything after it down).  After the loop, paragraphs 5..20
#    hold the license block text.
# ---------------------------------------------------------------------------
$lines = @(
    "Beyond compare 4 license",
    "=========================================================",
    "Licensed to:    Honeywell",
    "Quantity:       World-wide users",
    "Serial number:  4872-6125",
    "License type:   BC4 Standard Edition, Multi Platform",
    "=========================================================",
    "--- BEGIN LICENSE KEY ---",
    "0q3ze4J269BJSw1Iynr73noxKQ5GOWBTukdkTRC2XbwHEGkIuR7wxHJ10",
    "soxC1FvLDqNPUSHHwrYM-ZszG6DbLqNftd-n-gI46jlO8pjxHtqQQKEUH",
    "s26GeLHMtaluBjOuogOwi4SOXGOt1jQEQeo9fldA7ysvhKBZGjDuGr9iv",
    "0siiIRiFSqeqmbb2+mcy-N1Hq-vARl4nWNxPpzBxSg9HnapAB+LXenBtR",
    "SuLrYEkwgagDhgHrK4aQLWGHzgBFS5JaNVFGqTS9T7gwYx4TmWpIZjg9A",
    "n1tB5oYRV5M8bHgQU3E0MREbf79lwSs-f80-QybSkYf63P5RiNRFEOODk",
    "--- END LICENSE KEY -----",
    "========================================================="
)

for ($i = 0; $i -lt $lines.Length; $i++) {
    $p = $d.Paragraphs.Item(5 + $i)
    $p.Range.Text = $lines[$i]
}

# The license block paragraphs (all but the very first "Beyond compare 4
# license" line) carry an explicit single-line-spacing override.
for ($i = 1; $i -lt $lines.Length; $i++) {
    $p = $d.Paragraphs.Item(5 + $i)
    $p.Range.ParagraphFormat.LineSpacingRule = 0
}

# ---------------------------------------------------------------------------
# 4. Paragraph 21 is the (still empty) paragraph that used to be the first
#    placeholder; it loses its bullet and becomes the new home of the
#    "_GoBack" bookmark.
# ---------------------------------------------------------------------------
$tail = $d.Paragraphs.Item(21)
$tail.Range.ListFormat.RemoveNumbers()
$d.Bookmarks.Add("_GoBack", $tail.Range)
